# Apply updated water-mass gamma_n (neutral density) limit values for V2.0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 27.95
$ws.Range("C5").Value = 27.95
$ws.Range("D5").Value = 28
$ws.Range("C6").Value = 28
$ws.Range("C9").Value = 27.85
$ws.Range("D9").Value = 27.95
$ws.Range("C10").Value = 28
$ws.Range("D11").Value = 27.1
$ws.Range("D12").Value = 27.75
$ws.Range("C13").Value = 27.85
$ws.Range("D13").Value = 27.9
$ws.Range("C14").Value = 27.9
$ws.Range("D14").Value = 28
$ws.Range("C15").Value = 28
$ws.Range("C18").Value = 27.8
$ws.Range("D18").Value = 27.9
$ws.Range("C19").Value = 27.9
$ws.Range("D19").Value = 28
$ws.Range("C20").Value = 28
$ws.Range("D23").Value = 27.1
$ws.Range("C24").Value = 27.45
$ws.Range("D24").Value = 27.65
$ws.Range("C25").Value = 27.8
$ws.Range("D25").Value = 27.9
$ws.Range("C26").Value = 27.9
$ws.Range("D26").Value = 28
$ws.Range("C27").Value = 28
$ws.Range("D30").Value = 26.7
$ws.Range("C31").Value = 26.7
$ws.Range("D31").Value = 27
$ws.Range("D32").Value = 27.65
$ws.Range("C33").Value = 27.7
$ws.Range("D33").Value = 27.85
$ws.Range("C34").Value = 27.85
$ws.Range("C35").Value = 28
$ws.Range("D39").Value = 27
$ws.Range("D40").Value = 27.6
$ws.Range("C41").Value = 27.75
$ws.Range("D41").Value = 27.9
$ws.Range("C42").Value = 27.9
$ws.Range("D42").Value = 28
$ws.Range("C43").Value = 28
$ws.Range("D46").Value = 27
$ws.Range("C47").Value = 27.3
$ws.Range("D47").Value = 27.6
$ws.Range("C48").Value = 27.8
$ws.Range("D48").Value = 27.9
$ws.Range("C49").Value = 27.9
$ws.Range("D49").Value = 28
$ws.Range("C50").Value = 28
$ws.Range("D52").Value = 27
$ws.Range("D53").Value = 27.6
$ws.Range("C54").Value = 27.85
$ws.Range("D54").Value = 27.95
$ws.Range("C55").Value = 27.95
$ws.Range("D55").Value = 28.05
$ws.Range("C56").Value = 28.05
$ws.Range("D59").Value = 27
$ws.Range("C60").Value = 27.2
$ws.Range("D60").Value = 27.6
$ws.Range("C61").Value = 27.8
$ws.Range("D61").Value = 27.9
$ws.Range("C62").Value = 27.9
$ws.Range("D62").Value = 28
$ws.Range("C63").Value = 28
$ws.Range("D66").Value = 27
$ws.Range("D67").Value = 27.6
$ws.Range("C68").Value = 27.8
$ws.Range("D68").Value = 27.9
$ws.Range("C69").Value = 27.9
$ws.Range("D69").Value = 28
$ws.Range("C70").Value = 28
$ws.Range("D73").Value = 27.5
$ws.Range("C74").Value = 27.85
$ws.Range("D74").Value = 27.95
$ws.Range("C75").Value = 27.95
$ws.Range("D75").Value = 28.05
$ws.Range("C76").Value = 28.05

# Reposition the view/selection to match the saved window state
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
$ws.Range("F35").Select()
